$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlCenter constant for HorizontalAlignment
$xlCenter = -4108

# --- Column C ("x" mark for attendance) for existing rows 2-17 ---
# Rows that get an "x" value in column C
$xRows = 2,3,4,5,6,9,10,11,17
foreach ($r in $xRows) {
    $ws.Range("C$r").Value = "x"
    $ws.Range("C$r").HorizontalAlignment = $xlCenter
}

# Rows that get the style applied to column C but remain empty
$emptyStyledRows = 7,8,12,13,14,15,16
foreach ($r in $emptyStyledRows) {
    $ws.Range("C$r").HorizontalAlignment = $xlCenter
}

# --- New rows 18-20 with new names in column A and "x" in column C ---
$ws.Range("A18").Value = "Epsilon"
$ws.Range("C18").Value = "x"
$ws.Range("C18").HorizontalAlignment = $xlCenter

$ws.Range("A19").Value = "Esperanza"
$ws.Range("C19").Value = "x"
$ws.Range("C19").HorizontalAlignment = $xlCenter

$ws.Range("A20").Value = "Santana"
$ws.Range("C20").Value = "x"
$ws.Range("C20").HorizontalAlignment = $xlCenter

# --- Update selection to match the edited workbook's cursor position ---
$ws.Range("C21").Select()
